$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Nuove righe: aggiornamento dati fino al 6 gennaio 2022 (date seriali 44539-44566,
# righe foglio 465-491). Colonne: A=data, B=nuovi pos., C=somma mobile 7gg.,
# D=somma mobile 7gg. per 100mila abitanti.
$firstNewRow = 465
$lastNewRow = 491

# Estendi la formattazione (stile/numero data) della colonna A nelle nuove righe,
# copiando dalla riga precedente (ultima riga esistente, 464).
$fillRange = "A" + $firstNewRow + ":A" + $lastNewRow
$ws.Range("A464").Copy() | Out-Null
$ws.Range($fillRange).PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = 0

$newData = @(
    @(44539,64,275,145.4926380725135),
    @(44540,60,280,148.1379587647411),
    @(44541,9,247,130.6788421960394),
    @(44542,80,299,158.1901773952056),
    @(44543,69,313,165.5970753334427),
    @(44544,27,316,167.1842677487792),
    @(44545,11,320,169.3005243025612),
    @(44546,60,316,167.1842677487792),
    @(44547,73,329,174.0621015485707),
    @(44548,83,403,213.212847793538),
    @(44550,87,410,216.9162967626565),
    @(44551,111,452,239.1369905773677),
    @(44552,18,443,234.3754133313581),
    @(44553,148,580,306.8572002983922),
    @(44554,87,607,321.1419320364208),
    @(44555,107,641,339.1301127435679),
    @(44556,137,695,367.6995762196251),
    @(44557,177,785,415.3153486797204),
    @(44558,184,858,453.9370307862422),
    @(44559,259,1099,581.4414881516086),
    @(44560,397,1348,713.1784586245391),
    @(44561,149,1410,745.9804352081603),
    @(44562,195,1498,792.5380793913646),
    @(44563,228,1589,840.6829159899055),
    @(44564,345,1757,929.5656912487501),
    @(44565,297,1870,989.349938893092),
    @(44566,499,2110,1116.325332120013)
)

$r = $firstNewRow
foreach ($row in $newData) {
    $ws.Cells.Item($r, 1).Value = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
    $ws.Cells.Item($r, 4).Value = $row[3]
    $r++
}

Write-Output "Aggiornate righe $firstNewRow-$lastNewRow"
